$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 171, shifting old rows
# 171-191 down to 174-194.
$ws.Rows("171:173").Insert()

# New row 171: Cuatro cascos verde
$ws.Cells.Item(171, 1).Value = 11
$ws.Cells.Item(171, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(171, 3).Value = "Bíobío"
$ws.Cells.Item(171, 4).Value = 44491
$ws.Cells.Item(171, 5).Value = 8
$ws.Cells.Item(171, 6).Value = 100112002
$ws.Cells.Item(171, 7).Value = "Pimiento"
$ws.Cells.Item(171, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 100
$ws.Cells.Item(171, 11).Value = 32000
$ws.Cells.Item(171, 12).Value = 34000
$ws.Cells.Item(171, 13).Value = 33000
$ws.Cells.Item(171, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(171, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(171, 16).Value = 1833
$ws.Cells.Item(171, 17).Value = 18
$ws.Cells.Item(171, 18).Value = "Hortaliza"

# New row 172: Zafiro rojo
$ws.Cells.Item(172, 1).Value = 11
$ws.Cells.Item(172, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(172, 3).Value = "Bíobío"
$ws.Cells.Item(172, 4).Value = 44491
$ws.Cells.Item(172, 5).Value = 8
$ws.Cells.Item(172, 6).Value = 100112002
$ws.Cells.Item(172, 7).Value = "Pimiento"
$ws.Cells.Item(172, 8).Value = "Zafiro rojo"
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 100
$ws.Cells.Item(172, 11).Value = 44000
$ws.Cells.Item(172, 12).Value = 45000
$ws.Cells.Item(172, 13).Value = 44500
$ws.Cells.Item(172, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(172, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(172, 16).Value = 2967
$ws.Cells.Item(172, 17).Value = 15
$ws.Cells.Item(172, 18).Value = "Hortaliza"

# New row 173: Zafiro verde
$ws.Cells.Item(173, 1).Value = 11
$ws.Cells.Item(173, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(173, 3).Value = "Bíobío"
$ws.Cells.Item(173, 4).Value = 44491
$ws.Cells.Item(173, 5).Value = 8
$ws.Cells.Item(173, 6).Value = 100112002
$ws.Cells.Item(173, 7).Value = "Pimiento"
$ws.Cells.Item(173, 8).Value = "Zafiro verde"
$ws.Cells.Item(173, 9).Value = "Primera"
$ws.Cells.Item(173, 10).Value = 100
$ws.Cells.Item(173, 11).Value = 32000
$ws.Cells.Item(173, 12).Value = 33000
$ws.Cells.Item(173, 13).Value = 32500
$ws.Cells.Item(173, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(173, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(173, 16).Value = 2167
$ws.Cells.Item(173, 17).Value = 15
$ws.Cells.Item(173, 18).Value = "Hortaliza"
